$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.947.31"
$ws.Range("E2").Value = "  +6.02%  "
$ws.Range("D3").Value = "3.298.47"
$ws.Range("E3").Value = "  +1.30%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'409.14"
$ws.Range("E5").Value = "  +3.37%  "
$ws.Range("D6").Value = "'111.49"
$ws.Range("E6").Value = "  +2.90%  "
$ws.Range("D7").Value = "3.294.14"
$ws.Range("E7").Value = "  +1.33%  "
$ws.Range("D8").Value = "'0.565"
$ws.Range("E8").Value = "  -4.79%  "
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("D10").Value = "'0.617"
$ws.Range("E10").Value = "  -1.67%  "
$ws.Range("D11").Value = "'0.113"
$ws.Range("E11").Value = "  +13.70%  "
$ws.Range("D12").Value = "'38.47"
$ws.Range("E12").Value = "  -1.68%  "
$ws.Range("E13").Value = "  -0.15%  "
$ws.Range("D14").Value = "3.822.27"
$ws.Range("E14").Value = "  +1.51%  "
$ws.Range("E15").Value = "  -1.19%  "
$ws.Range("D16").Value = "'18.87"
$ws.Range("E16").Value = "  -1.31%  "
$ws.Range("D17").Value = "3.309.06"
$ws.Range("E17").Value = "  +1.19%  "
$ws.Range("D18").Value = "60.769.41"
$ws.Range("E18").Value = "  +6.26%  "
$ws.Range("D19").Value = "'0.981"
$ws.Range("E19").Value = "  -4.86%  "
$ws.Range("D20").Value = "'10.44"
$ws.Range("E20").Value = "  -3.01%  "
$ws.Range("E21").Value = "  +0.27%  "
$ws.Range("E22").Value = "  -3.39%  "
$ws.Range("D23").Value = "'12.50"
$ws.Range("E23").Value = "  -3.22%  "
$ws.Range("D24").Value = "'295.37"
$ws.Range("E24").Value = "  -0.36%  "
$ws.Range("D25").Value = "'72.99"
$ws.Range("E25").Value = "  -1.96%  "
$ws.Range("E26").Value = "  -2.74%  "
$ws.Range("D27").Value = "'29.07"
$ws.Range("E27").Value = "  +3.89%  "
$ws.Range("D28").Value = "'4.27"
$ws.Range("E28").Value = "  -2.07%  "
$ws.Range("D29").Value = "'0.172"
$ws.Range("E29").Value = "  +2.82%  "
$ws.Range("D30").Value = "'7.30"
$ws.Range("E30").Value = "  +0.49%  "
$ws.Range("D31").Value = "'7.42"
$ws.Range("E31").Value = "  -2.31%  "
$ws.Range("D32").Value = "'1.00"
$ws.Range("E32").Value = "  +0.04%  "
$ws.Range("E33").Value = "  +2.03%  "
$ws.Range("D34").Value = "'11.07"
$ws.Range("E34").Value = "  -3.13%  "
$ws.Range("E35").Value = "  +15.14%  "
$ws.Range("D36").Value = "'38.57"
$ws.Range("E36").Value = "  -1.42%  "
$ws.Range("E37").Value = "  -0.76%  "
$ws.Range("D38").Value = "'52.21"
$ws.Range("E38").Value = "  +1.16%  "
$ws.Range("D39").Value = "'0.998"
$ws.Range("E39").Value = "  +0.15%  "
$ws.Range("D40").Value = "'3.05"
$ws.Range("E40").Value = "  +4.32%  "
$ws.Range("E41").Value = "  -6.06%  "
$ws.Range("D42").Value = "'134.25"
$ws.Range("E42").Value = "  -0.08%  "
$ws.Range("B43").Value = "TheGraph"
$ws.Range("C43").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D43").Value = "'0.286"
$ws.Range("E43").Value = "  +1.82%  "
$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D44").Value = "'0.119"
$ws.Range("E44").Value = "  -3.11%  "
$ws.Range("E45").Value = "  -0.91%  "
$ws.Range("D46").Value = "'16.13"
$ws.Range("E46").Value = "  -5.89%  "
$ws.Range("E47").Value = "  -5.60%  "
$ws.Range("E48").Value = "  +2.04%  "
$ws.Range("D49").Value = "'20.82"
$ws.Range("E49").Value = "  -5.92%  "
$ws.Range("D50").Value = "2.100.15"
$ws.Range("E50").Value = "  -2.79%  "
$ws.Range("D51").Value = "3.631.47"
$ws.Range("E51").Value = "  +1.66%  "
